$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.779.97'
$ws.Range('E2').Value = '  +0.25%  '
$ws.Range('D3').Value = '3.548.63'
$ws.Range('E3').Value = '  +3.28%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '599.73'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +3.08%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '135.18'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +0.46%  '
$ws.Range('D7').Value = '3.547.20'
$ws.Range('E7').Value = '  +3.22%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.493'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +1.89%  '
$ws.Range('E10').Value = '  +1.21%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '6.89'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -1.52%  '
$ws.Range('E12').Value = '  +2.31%  '
$ws.Range('D13').Value = '4.152.74'
$ws.Range('E13').Value = '  +3.39%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.0000181'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +1.32%  '
$ws.Range('D15').Value = '3.555.17'
$ws.Range('E15').Value = '  +3.46%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '26.92'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +2.65%  '
$ws.Range('E17').Value = '  +0.79%  '
$ws.Range('D18').Value = '64.679.23'
$ws.Range('E18').Value = '  +0.30%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '10.00'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +3.52%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '14.32'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +5.35%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '5.81'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +2.11%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '386.10'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +1.25%  '
$ws.Range('E23').Value = '  +5.73%  '
$ws.Range('D24').Value = '3.692.50'
$ws.Range('E24').Value = '  +3.44%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '74.22'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +3.07%  '
$ws.Range('E26').Value = '  +0.07%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.0000116'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +9.95%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '7.60'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +6.51%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.999'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +0.03%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '2.29'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +4.94%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '8.34'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +3.35%  '
$ws.Range('D32').Value = '3.559.36'
$ws.Range('E32').Value = '  +3.25%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.45'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +23.24%  '
$ws.Range('B34').Value = 'EthereumClassic'
$ws.Range('C34').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '23.93'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +4.83%  '
$ws.Range('B35').Value = 'USDe'
$ws.Range('C35').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.00'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +0.00%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.143'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +2.33%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '169.49'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -0.42%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '6.90'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +2.48%  '
$ws.Range('E39').Value = '  +5.53%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '4.97'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +7.11%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.0801'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +5.15%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.826'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +3.12%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '26.82'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +18.72%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '42.56'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +1.61%  '
$ws.Range('E45').Value = '  +0.19%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '4.44'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +3.82%  '
$ws.Range('E47').Value = '  +9.50%  '
$ws.Range('E48').Value = '  +3.50%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '6.91'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +5.73%  '
$ws.Range('D50').Value = '2.449.69'
$ws.Range('E50').Value = '  +11.48%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.36'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +14.23%  '
